$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and report week date range)
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# Crime statistics table updates (rows 14-30)
$ws.Range("M14").Value = -73.333333333333
$ws.Range("N14").Value = -81.818181818181
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -11.111111111111
$ws.Range("M15").Value = -11.111111111111
$ws.Range("C16").Value = 3
$ws.Range("I29").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 75
$ws.Range("K16").Value = 47.058823529411
$ws.Range("L16").Value = -15.730337078651
$ws.Range("M16").Value = -76.038338658147
$ws.Range("N16").Value = -92.096944151738
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -35.294117647058
$ws.Range("I17").Value = 208
$ws.Range("J17").Value = 239
$ws.Range("K17").Value = -12.970711297071
$ws.Range("L17").Value = -6.726457399103
$ws.Range("M17").Value = -20.610687022900
$ws.Range("N17").Value = -46.113989637305
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -70.588235294117
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -9.243697478991
$ws.Range("L18").Value = -28
$ws.Range("M18").Value = -69.91643454039
$ws.Range("N18").Value = -92.833443928334
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -53.333333333333
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -17.142857142857
$ws.Range("I19").Value = 380
$ws.Range("J19").Value = 328
$ws.Range("K19").Value = 15.853658536585
$ws.Range("L19").Value = -3.553299492385
$ws.Range("M19").Value = -8.212560386473
$ws.Range("N19").Value = -30.783242258652
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -47.619047619047
$ws.Range("J20").Value = 256
$ws.Range("K20").Value = -32.421875
$ws.Range("L20").Value = 19.310344827586
$ws.Range("M20").Value = -51.404494382022
$ws.Range("N20").Value = -94.573400250941
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -50
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -34.736842105263
$ws.Range("I21").Value = 964
$ws.Range("J21").Value = 1011
$ws.Range("K21").Value = -4.648862512364
$ws.Range("L21").Value = -5.024630541871
$ws.Range("M21").Value = -44.502014968336
$ws.Range("N21").Value = -85.475365375922
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 150
$ws.Range("F24").Value = 62
$ws.Range("G24").Value = 36
$ws.Range("H24").Value = 72.222222222222
$ws.Range("I24").Value = 647
$ws.Range("J24").Value = 613
$ws.Range("K24").Value = 5.546492659053
$ws.Range("L24").Value = -11.852861035422
$ws.Range("M24").Value = -14.531043593130
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 152
$ws.Range("J25").Value = 153
$ws.Range("K25").Value = -0.653594771241
$ws.Range("L25").Value = 10.948905109489
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 26
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -7.142857142857
$ws.Range("I26").Value = 354
$ws.Range("J26").Value = 374
$ws.Range("K26").Value = -5.347593582887
$ws.Range("L26").Value = 3.206997084548
$ws.Range("M26").Value = -36.101083032491
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = -34.615384615384
$ws.Range("L27").Value = -22.727272727272
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("J28").Value = 37
$ws.Range("K28").Value = -5.405405405405
$ws.Range("M29").Value = -97.872340425531
$ws.Range("N29").Value = -98.701298701298
$ws.Range("M30").Value = -96.774193548387
$ws.Range("N30").Value = -98.529411764705
